$d = $word.ActiveDocument

# --- 1) Replace the TODO block (SamplerImpl.../FailureCounter/.../SamplerImplStressTest/System-wide Metrics section)
#        with the new set of bullet paragraphs described by the commit. ---

$startRange = $d.Content
$null = $startRange.Find.Execute("(go to SamplerImpl line", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$blockStart = $startRange.Paragraphs.Item(1).Range.Start

$endRange = $d.Content
$null = $endRange.Find.Execute("Future parameters", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$blockEnd = $endRange.Paragraphs.Item(1).Range.End

$target = $d.Range($blockStart, $blockEnd)

$newParas = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080"/><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:b/></w:rPr><w:t>A start()/record()/stop() must generate at least a sample. Test for that.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080"/><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t>Apply this when implementing metrics, below:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:br/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:b/></w:rPr><w:t>S</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:b/></w:rPr><w:t>ystem-wide</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:b/></w:rPr><w:t xml:space="preserve"> Metrics</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:b/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t>refactor – I should register types with the sampler and get values. Now I am registering values and getting values – not intuitive.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t xml:space="preserve"> Also every time I invoke getValue() I get the current value, not the one recorded at the time of the snapshot.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:br/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:b/></w:rPr><w:t>SampleInterval extrapolation must also extrapolate metrics.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:b/></w:rPr><w:br/></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t>Return to SamplerTest.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t>registerAndReadMetrics</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t>()</w:t></w:r></w:p>'

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newParas + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($packageXml)

# --- 2) Drop the stale lastRenderedPageBreak cached before the "log4j" run. ---

$lr = $d.Content
$null = $lr.Find.Execute("log4j. Externalize", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$log4jPara = $lr.Paragraphs.Item(1).Range

$log4jInner = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/><w:b/></w:rPr><w:t>log4j</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr><w:t>. Externalize log4j.xml so I can modify it directly, don't have to unjar for it. Ideally, the debugging should be turned on by --debug.</w:t></w:r></w:p>
'@

$log4jPackage = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $log4jInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$log4jPara.InsertXML($log4jPackage)

Write-Output "Done."
